$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "KETOPREK 75 MG 20 CAPS." line item (row 35) by pulling every row below it
# up by one, preserving each row's own formatting (row height stays tied to its row
# number rather than travelling with the shifted-up data).
for ($r = 35; $r -le 77; $r++) {
    $srcRow = $r + 1
    for ($col = 1; $col -le 14; $col++) {
        $src = $ws.Cells.Item($srcRow, $col)
        $dst = $ws.Cells.Item($r, $col)
        $dst.Value = $src.Value2
    }
}

# The totals row (previously row 77, now row 76) holds a hard-coded sum that must be
# reduced by the deleted item's price (18) to stay consistent with the remaining rows.
$ws.Range("K76").Value = 5638.73

# Drop the now-duplicated trailing row (what used to be row 78) so the sheet is one
# row shorter overall, matching a genuine row deletion.
$ws.Rows("78:78").Delete()
